$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '54.324.14'
$ws.Range("E2").Value = '  +0.82%  '

$ws.Range("D3").Value = '2.242.37'
$ws.Range("E3").Value = '  +0.14%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.15%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '494.96'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.51%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '127.32'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.56%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.993'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.54%  '

$ws.Range("E8").Value = '  +1.44%  '

$ws.Range("D9").Value = '2.284.57'
$ws.Range("E9").Value = '  +1.46%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0946'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.86%  '

$ws.Range("E11").Value = '  +2.44%  '

$ws.Range("E12").Value = '  +3.00%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.62'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.24%  '

$ws.Range("D14").Value = '2.650.08'
$ws.Range("E14").Value = '  +0.45%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '21.71'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.64%  '

$ws.Range("D16").Value = '54.184.11'
$ws.Range("E16").Value = '  +0.72%  '

$ws.Range("E17").Value = '  +0.78%  '

$ws.Range("D18").Value = '2.262.81'
$ws.Range("E18").Value = '  +1.87%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.02'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.40%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.09'
$ws.Range("D20").Style = "Normal"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.47'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +6.03%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '301.30'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.41%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.995'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.57%  '

$ws.Range("E24").Value = '  -2.23%  '

$ws.Range("E25").Value = '  -2.18%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.996'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.41%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.372'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.89%  '

$ws.Range("D28").Value = '2.377.80'
$ws.Range("E28").Value = '  +0.75%  '

$ws.Range("E29").Value = '  +4.52%  '

$ws.Range("E30").Value = '  +0.27%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '168.42'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.38%  '

$ws.Range("D32").Value = '0.0₃0688'
$ws.Range("E32").Value = '  -0.07%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.60'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.14%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.86'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.88%  '

$ws.Range("E35").Value = '  -0.25%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.993'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.50%  '

$ws.Range("E37").Value = '  +0.75%  '

$ws.Range("E38").Value = '  +0.75%  '

$ws.Range("E39").Value = '  +2.31%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.861'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.86%  '

$ws.Range("E41").Value = '  +2.94%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '35.35'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.92%  '

$ws.Range("B43").Value = 'Stacks'
$ws.Range("C43").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.40'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.17%  '

$ws.Range("B44").Value = 'PolygonEcosystemToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.374'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.92%  '

$ws.Range("E45").Value = '  +1.16%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '128.11'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.60%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.92'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +5.30%  '

$ws.Range("E48").Value = '  +1.08%  '

$ws.Range("E49").Value = '  +0.96%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '237.55'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.34%  '

$ws.Range("E51").Value = '  +2.36%  '
